$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F9").Value = 'client_number=1010, first_name=" ", last_name="Warya", email_address="sukhtabwarya@gmail.com .com"'
$ws.Range("G9").Value = 'ValueError with message "First name cannot be blank."'

$ws.Range("F10").Value = 'client_number=1010, first_name="Sukhtab", last_name=" ", email_address="sukhtabwarya@gmail.com"'
$ws.Range("G10").Value = 'ValueError with message "Last name cannot be blank."'

$ws.Range("F11").Value = 'client_number=1010, first_name="Sukhtab", last_name="Warya", email_address="invalid-email"'
$ws.Range("G11").Value = 'email_address should be "email@pixell-river.com"'

$ws.Range("F12").Value = "Client instance created with client_number=1010"
$ws.Range("G12").Value = "client_number should be 1010"

$ws.Range("F13").Value = 'Client instance created with first_name="Sukhtab"'
$ws.Range("G13").Value = 'first_name should be "Sukhtab"'

$ws.Range("F14").Value = 'Client instance created with last_name="Warya"'
$ws.Range("G14").Value = 'last_name should be "Warya"'

$ws.Range("F15").Value = 'Client instance created with email_address="sukhtabwarya@gmail.com"'
$ws.Range("G15").Value = 'email_address should be "sukhtabwarya@gmail.com"'

$ws.Range("F16").Value = 'Client instance created with client_number=1010, first_name="Sukhtab", last_name="Warya", email_address="sukhtabwarya@gmail.com"'
$ws.Range("G16").Value = 'The string should be "Warya, Sukhtab[1010] - sukhtabwarya@gmail.com\n"'

$ws.Range("G7").Value = 'client_number should be 1010, first_name should be "Sukhtab", last_name should be "Warya", email_address should be "sukhtabwarya@gmail.com"'

# Update the selected cell / top-left view to reflect the saved view state
$ws.Range("F18").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 3
